# Apply the Tottenham_stats.xlsx edit:
#  1. Rename the stat sheets to human-friendly (spaced) names.
#  2. Bump the "Age" column (format "YY-DDD") day-of-year counter by 1 for
#     every player row on every stat sheet (data refreshed one day later).
#  3. On "Standard Stats" and "Playing Time" the "Playing Time" merged
#     header group shifted one column right (F1:I1 -> G1:I1), with F1
#     turning into a blank "Unnamed: 4_level_0" pandas-export placeholder
#     column header.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1) Rename sheets (index is stable; only the label text changes)
# ---------------------------------------------------------------------
$renames = @{
    2  = "Standard Stats"
    3  = "Shooting Stats"
    4  = "Passing Stats"
    5  = "Pass Types"
    6  = "Goal & Shot Creation"
    7  = "Defensive Actions"
    9  = "Playing Time"
    10 = "Miscellaneous Stats"
}
foreach ($idx in $renames.Keys) {
    $wb.Worksheets.Item($idx).Name = $renames[$idx]
}

# ---------------------------------------------------------------------
# 2) Age column bump, rows 4..lastRow, column E, on every stat sheet
# ---------------------------------------------------------------------
$lastRows = @{
    2  = 43   # Standard Stats
    3  = 34   # Shooting Stats
    4  = 34   # Passing Stats
    5  = 34   # Pass Types
    6  = 34   # Goal & Shot Creation
    7  = 34   # Defensive Actions
    8  = 34   # Possession
    9  = 43   # Playing Time
    10 = 34   # Miscellaneous Stats
}

foreach ($idx in $lastRows.Keys) {
    $ws = $wb.Worksheets.Item($idx)
    $last = $lastRows[$idx]
    for ($r = 4; $r -le $last; $r++) {
        $cell = $ws.Range("E$r")
        $txt = [string]$cell.Text
        if ($txt -match "^(\d+)-(\d+)$") {
            $years = $matches[1]
            $days = [int]$matches[2]
            $newDays = $days + 1
            $newTxt = "{0}-{1:D3}" -f $years, $newDays
            $cell.Value = $newTxt
        }
    }
}

# ---------------------------------------------------------------------
# 3) "Playing Time" merged header: F1:I1 -> G1:I1 on sheets 2 and 9
# ---------------------------------------------------------------------
foreach ($idx in 2, 9) {
    $ws = $wb.Worksheets.Item($idx)

    $ws.Range("F1:I1").UnMerge()
    $ws.Range("G1:I1").Merge()

    # Carry the header formatting (border/font) from F1 onto the newly
    # freed G1:I1 cells so they keep looking like the rest of the band.
    $ws.Range("F1").Copy()
    $ws.Range("G1:I1").PasteSpecial(-4122)   # xlPasteFormats

    $ws.Range("G1").Value = "Playing Time"
    $ws.Range("F1").Value = "Unnamed: 4_level_0"
}
